# Lecture partielle de l'EDT M1 MIAGE.
# Shift every week's date forward by 1096 days (same weekday family moves
# from "dimanche" to "jeudi"), and fix the last session's start time back
# to 13:30 (it incorrectly read 14:30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows holding the week-start date in column A / weekday label in column B.
$dateRows = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29, 32)

foreach ($r in $dateRows) {
    $oldSerial = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $oldSerial + 1096
}

# The weekday label column (B) always read "dimanche"; the schedule now
# starts on "jeudi" instead, so replace every occurrence in the sheet.
$ws.Cells.Replace("dimanche", "jeudi")

# Row 33 (last course block) had a stray "14:30" start time; correct it to
# match every other week's "13:30".
$ws.Range("D33").Value = "13:30"
